$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraphs we need to touch by their current text content so the
# script isn't brittle to exact paragraph-index assumptions.
$pStart = $null   # "Minuets:"
$pStop = $null    # "Desc:"  (has spell-check proofErr markup around "Desc")
$pDesc = $null    # empty paragraph holding the _GoBack bookmark
$pLast = $null    # final (last) paragraph in the document

foreach ($p in $d.Paragraphs) {
    # Range.Text includes the trailing paragraph-mark character(s); strip them
    # before comparing so exact-text matches work.
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Minuets:") { $pStart = $p }
    elseif ($t -eq "Desc:") { $pStop = $p }
    $pLast = $p
}
# The bookmark paragraph is the one right after the "Desc:" paragraph.
if ($pStop -ne $null) {
    $pDesc = $pStop.Next()
}

# 1) "Minuets:" paragraph -> "Start" + ":" + " 8:00 AM" (three separate runs)
$xmlStart = '<w:p ' + $wNs + '><w:r><w:t>Start</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> 8:00 AM</w:t></w:r></w:p>'
$pStart.Range.InsertXML($xmlStart)

# 2) "Desc:" paragraph (with spell-check proof markup) -> "Stop:" (single run, no proofErr)
$xmlStop = '<w:p ' + $wNs + '><w:r><w:t>Stop:</w:t></w:r></w:p>'
$pStop.Range.InsertXML($xmlStop)

# 3) Bookmark paragraph -> "Desc" + "ription" runs, bookmark preserved in the middle, then ":" run
$xmlDesc = '<w:p ' + $wNs + '><w:r><w:t>Desc</w:t></w:r><w:r><w:t>ription</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>:</w:t></w:r></w:p>'
$pDesc.Range.InsertXML($xmlDesc)

# 4) Insert an extra empty paragraph right before the final trailing empty paragraph
$insAt = $d.Range($pLast.Range.Start, $pLast.Range.Start)
$xmlEmpty = '<w:p ' + $wNs + '/>'
$insAt.InsertXML($xmlEmpty)
